# Add a new worksheet "dummy" at the end of the workbook, as a copy of the
# last existing efficiency-table sheet (so it inherits the same header style,
# borders, sheetPr/pageMargins layout), then overwrite its data with the new
# trimming-efficiency dataset (31 data rows, A1:E32).
$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$templateSheet = $wb.Worksheets.Item($sheetCount)
$lastSheet = $wb.Worksheets.Item($sheetCount)

# Copy-after-last duplicates the template sheet (formatting/styles included)
# and places the duplicate as the new last sheet.
$templateSheet.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "dummy"

# Header row (A1:E1) is inherited from the template (vbat_V / ibat_A / vbus_V /
# ibus_A / efficiency, bold+bordered+centered style) -- no change needed there.

$data = New-Object 'object[,]' 31,5
$data[0,0] = 3.9890679875
$data[0,1] = 0.1004135
$data[0,2] = 9.004001000000001
$data[0,3] = 0.05225086
$data[0,4] = 85.14018035526433
$data[1,0] = 3.9890748275
$data[1,1] = 0.100416
$data[1,2] = 9.004002
$data[1,3] = 0.05227101
$data[1,4] = 85.10961499425376
$data[2,0] = 3.9902983525
$data[2,1] = 0.2004981
$data[2,2] = 9.003961
$data[2,3] = 0.09682997
$data[2,4] = 91.76397708382639
$data[3,0] = 3.9918559025
$data[3,1] = 0.3004653
$data[3,2] = 9.003992999999999
$data[3,3] = 0.1416674
$data[3,4] = 94.02949578439697
$data[4,0] = 3.9936480175
$data[4,1] = 0.4004739
$data[4,2] = 9.003996000000001
$data[4,3] = 0.1866496
$data[4,4] = 95.1659627777638
$data[5,0] = 3.9950499225
$data[5,1] = 0.5004818
$data[5,2] = 9.003871
$data[5,3] = 0.2320144
$data[5,4] = 95.71197881913486
$data[6,0] = 3.996719685
$data[6,1] = 0.6005057
$data[6,2] = 9.003921
$data[6,3] = 0.2776032
$data[6,4] = 96.02065844021489
$data[7,0] = 3.998580715
$data[7,1] = 0.7005562
$data[7,2] = 9.003899000000001
$data[7,3] = 0.3234304
$data[7,4] = 96.1916546736407
$data[8,0] = 4.000355595
$data[8,1] = 0.8005186
$data[8,2] = 9.004004
$data[8,3] = 0.3694859
$data[8,4] = 96.25792065245291
$data[9,0] = 4.00187433
$data[9,1] = 0.9005030000000001
$data[9,2] = 9.003869
$data[9,3] = 0.4155703
$data[9,4] = 96.3107890519925
$data[10,0] = 4.0030724275
$data[10,1] = 1.000478
$data[10,2] = 9.003776999999999
$data[10,3] = 0.4620328
$data[10,4] = 96.27276683247364
$data[11,0] = 4.00497874
$data[11,1] = 1.100541
$data[11,2] = 9.004151
$data[11,3] = 0.5087447
$data[11,4] = 96.21965030670587
$data[12,0] = 4.00678753
$data[12,1] = 1.200615
$data[12,2] = 9.004236000000001
$data[12,3] = 0.5556563
$data[12,4] = 96.14948589439399
$data[13,0] = 4.0081607925
$data[13,1] = 1.300532
$data[13,2] = 9.004201999999999
$data[13,3] = 0.6025578
$data[13,4] = 96.07761993095365
$data[14,0] = 4.00987909
$data[14,1] = 1.400625
$data[14,2] = 9.004047999999999
$data[14,3] = 0.6498704
$data[14,4] = 95.98173441924138
$data[15,0] = 4.011741625
$data[15,1] = 1.500519
$data[15,2] = 9.00433
$data[15,3] = 0.6973712
$data[15,4] = 95.86477175004329
$data[16,0] = 4.0132960225
$data[16,1] = 1.600656
$data[16,2] = 9.004076
$data[16,3] = 0.7452010999999999
$data[16,4] = 95.73848752412417
$data[17,0] = 4.0152297375
$data[17,1] = 1.700475
$data[17,2] = 9.004213999999999
$data[17,3] = 0.7931827
$data[17,4] = 95.60081817487043
$data[18,0] = 4.016991787499999
$data[18,1] = 1.800623
$data[18,2] = 9.004089
$data[18,3] = 0.841364
$data[18,4] = 95.47727873176729
$data[19,0] = 4.0184857
$data[19,1] = 1.900542
$data[19,2] = 9.003947
$data[19,3] = 0.8897482
$data[19,4] = 95.33225163876111
$data[20,0] = 4.0200871325
$data[20,1] = 2.000669
$data[20,2] = 9.004134000000001
$data[20,3] = 0.9384987
$data[20,4] = 95.17767335401078
$data[21,0] = 4.0217745675
$data[21,1] = 2.100584
$data[21,2] = 9.004168
$data[21,3] = 0.9874233
$data[21,4] = 95.01907891187091
$data[22,0] = 4.02373624
$data[22,1] = 2.200623
$data[22,2] = 9.004177
$data[22,3] = 1.036509
$data[22,4] = 94.87636806853352
$data[23,0] = 4.025525235
$data[23,1] = 2.300559
$data[23,2] = 9.004118999999999
$data[23,3] = 1.085759
$data[23,4] = 94.72863188166613
$data[24,0] = 4.02685571
$data[24,1] = 2.400647
$data[24,2] = 9.004130999999999
$data[24,3] = 1.135508
$data[24,4] = 94.55018209429991
$data[25,0] = 4.028533175
$data[25,1] = 2.500538
$data[25,2] = 9.004217000000001
$data[25,3] = 1.185483
$data[25,4] = 94.37112228487878
$data[26,0] = 4.03024146
$data[26,1] = 2.600561
$data[26,2] = 9.004092
$data[26,3] = 1.23532
$data[26,4] = 94.22772701581142
$data[27,0] = 4.0317491325
$data[27,1] = 2.700638
$data[27,2] = 9.004217000000001
$data[27,3] = 1.285704
$data[27,4] = 94.05305949086461
$data[28,0] = 4.0334121
$data[28,1] = 2.800624
$data[28,2] = 9.004341999999999
$data[28,3] = 1.336512
$data[28,4] = 93.8647566743014
$data[29,0] = 4.03513508
$data[29,1] = 2.900586
$data[29,2] = 9.004186000000001
$data[29,3] = 1.387176
$data[29,4] = 93.7060864205001
$data[30,0] = 4.03677894
$data[30,1] = 3.000629
$data[30,2] = 9.004253
$data[30,3] = 1.438235
$data[30,4] = 93.53404733163352

$ws.Range("A2:E32").Value2 = $data

Write-Output ("sheets=" + $wb.Worksheets.Count)
Write-Output ("name=" + $ws.Name)
